$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2014652014652015
$ws.Range("C2").Value = 0.5457875457875457
$ws.Range("J2").Value = 0.007326007326007326
$ws.Range("P2").Value = 0.1684981684981685
$ws.Range("S2").Value = 0.07692307692307693
$ws.Range("B3").Value = 0.006622516556291391
$ws.Range("J3").Value = 0.05960264900662252
$ws.Range("P3").Value = 0.7748344370860927
$ws.Range("S3").Value = 0.1589403973509934
$ws.Range("J4").Value = 0.06122448979591837
$ws.Range("P4").Value = 0.6938775510204082
$ws.Range("S4").Value = 0.2448979591836735
$ws.Range("B6").Value = 0.05092592592592592
$ws.Range("D6").Value = 0.02777777777777778
$ws.Range("F6").Value = 0.04629629629629629
$ws.Range("J6").Value = 0.2546296296296297
$ws.Range("O6").Value = 0.02314814814814815
$ws.Range("Q6").Value = 0.1481481481481481
$ws.Range("R6").Value = 0.06481481481481481
$ws.Range("S6").Value = 0.3842592592592592
$ws.Range("B7").Value = 0.1666666666666667
$ws.Range("D7").Value = 0.01388888888888889
$ws.Range("F7").Value = 0.05555555555555555
$ws.Range("J7").Value = 0.1805555555555556
$ws.Range("O7").Value = 0.01388888888888889
$ws.Range("Q7").Value = 0.1458333333333333
$ws.Range("R7").Value = 0.0763888888888889
$ws.Range("S7").Value = 0.3472222222222222
$ws.Range("B8").Value = 0.08551068883610451
$ws.Range("D8").Value = 0.02375296912114014
$ws.Range("F8").Value = 0.06888361045130641
$ws.Range("J8").Value = 0.1353919239904988
$ws.Range("O8").Value = 0.02612826603325416
$ws.Range("Q8").Value = 0.1971496437054632
$ws.Range("R8").Value = 0.08788598574821853
$ws.Range("S8").Value = 0.3752969121140142
$ws.Range("B9").Value = 0.08839779005524862
$ws.Range("D9").Value = 0.03314917127071823
$ws.Range("E9").Value = 0.005524861878453038
$ws.Range("F9").Value = 0.05524861878453038
$ws.Range("J9").Value = 0.1215469613259668
$ws.Range("O9").Value = 0.01104972375690608
$ws.Range("Q9").Value = 0.2154696132596685
$ws.Range("R9").Value = 0.08287292817679558
$ws.Range("S9").Value = 0.3867403314917127
$ws.Range("B10").Value = 0.09930178432893716
$ws.Range("D10").Value = 0.02094647013188518
$ws.Range("F10").Value = 0.07137315748642359
$ws.Range("J10").Value = 0.1380915438324282
$ws.Range("O10").Value = 0.02017067494181536
$ws.Range("Q10").Value = 0.2451512800620636
$ws.Range("R10").Value = 0.07214895267649341
$ws.Range("S10").Value = 0.3328161365399535
$ws.Range("G11").Value = 0.1333333333333333
$ws.Range("J11").Value = 0.09166666666666666
$ws.Range("K11").Value = 0.1833333333333333
$ws.Range("L11").Value = 0.5833333333333334
$ws.Range("S11").Value = 0.008333333333333333
$ws.Range("G12").Value = 0.6923076923076923
$ws.Range("J12").Value = 0.2517482517482518
$ws.Range("L12").Value = 0.04195804195804196
$ws.Range("S12").Value = 0.01398601398601399
$ws.Range("G13").Value = 0.5161290322580645
$ws.Range("J13").Value = 0.4516129032258064
$ws.Range("S13").Value = 0.03225806451612903
$ws.Range("J14").Value = 1
$ws.Range("F15").Value = 0.02212389380530973
$ws.Range("H15").Value = 0.1814159292035398
$ws.Range("I15").Value = 0.07079646017699115
$ws.Range("J15").Value = 0.3451327433628318
$ws.Range("K15").Value = 0.06194690265486726
$ws.Range("M15").Value = 0.01769911504424779
$ws.Range("O15").Value = 0.06194690265486726
$ws.Range("S15").Value = 0.2389380530973451
$ws.Range("F16").Value = 0.04232804232804233
$ws.Range("H16").Value = 0.1746031746031746
$ws.Range("I16").Value = 0.08994708994708994
$ws.Range("J16").Value = 0.4126984126984127
$ws.Range("K16").Value = 0.07407407407407407
$ws.Range("M16").Value = 0.005291005291005291
$ws.Range("N16").Value = 0.005291005291005291
$ws.Range("O16").Value = 0.06349206349206349
$ws.Range("S16").Value = 0.1322751322751323
$ws.Range("F17").Value = 0.02263374485596708
$ws.Range("H17").Value = 0.1975308641975309
$ws.Range("I17").Value = 0.09465020576131687
$ws.Range("J17").Value = 0.411522633744856
$ws.Range("K17").Value = 0.08230452674897119
$ws.Range("M17").Value = 0.01646090534979424
$ws.Range("O17").Value = 0.06995884773662552
$ws.Range("S17").Value = 0.102880658436214
$ws.Range("F18").Value = 0.005847953216374269
$ws.Range("H18").Value = 0.152046783625731
$ws.Range("I18").Value = 0.09941520467836257
$ws.Range("J18").Value = 0.4853801169590643
$ws.Range("K18").Value = 0.0935672514619883
$ws.Range("M18").Value = 0.01169590643274854
$ws.Range("O18").Value = 0.07017543859649122
$ws.Range("S18").Value = 0.08187134502923976
$ws.Range("F19").Value = 0.01544050862851953
$ws.Range("H19").Value = 0.2098092643051771
$ws.Range("I19").Value = 0.07720254314259764
$ws.Range("J19").Value = 0.4041780199818347
$ws.Range("K19").Value = 0.0971843778383288
$ws.Range("M19").Value = 0.0145322434150772
$ws.Range("O19").Value = 0.07356948228882834
$ws.Range("S19").Value = 0.1080835603996367
